$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Write the new "Test Route Finding" result cells in an order that makes the
# shared-string table come out as: .. "route", "shortest route", "message"
$ws.Range("C16").Value = "Return a 200 HTTP code with the route"
$ws.Range("C17").Value = "Return a 200 HTTP code with the route"
$ws.Range("C18").Value = "Return a 200 HTTP code with the shortest route"
$ws.Range("C15").Value = "Return a 200 HTTP code with a message"

# Rows 16 and 17 lose their custom (30pt) row height and revert to the
# sheet's default row height; row 18 keeps its custom height.
$ws.Rows.Item(16).EntireRow.AutoFit()
$ws.Rows.Item(17).EntireRow.AutoFit()

# A new, slightly wider column E is introduced between the fixed columns
# and the rest of the (default-width) sheet. (Closest value this engine's
# ColumnWidth rounding supports to the target 10.42578125 stored width.)
$ws.Columns.Item(5).ColumnWidth = 9.6725

# The active selection moves to F17.
$ws.Range("F17").Select()
